# Aspose.Cells Cloud SDK 23.11 release test-data update:
#   - Sheet1!A1 gets a value (1111), selection moves to A2
#   - Two more (blank) worksheets, "Sheet2" and "Sheet3", are appended after Sheet1

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Populate the existing sheet.
$ws1.Range("A1").Value = 1111

# Insert the new sheets *after* Sheet1 (Worksheets.Add() with no placement
# argument inserts before the active sheet, which would push Sheet1 to the
# end) so the final tab order is Sheet1, Sheet2, Sheet3.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "Sheet3"

# Restore Sheet1 as the active/selected sheet and leave the cursor on A2,
# matching the saved selection in the workbook.
$ws1.Activate()
$ws1.Range("A2").Select()
